# Applies the Barkley Sockeye stock-recruit infilled data update to the "S-R data" sheet.
# Updates historical HED (Hucuktlis) rows 96-129 with revised S/H/N/H_cv values,
# replaces rows 143-147 with an extended historical HED series (1918-1922),
# and appends new historical HED rows through 1976 (rows 148-201).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S-R data")

# Row 96
$ws.Cells.Item(96, 3).Value = 4800
$ws.Cells.Item(96, 8).Value = 30860
$ws.Cells.Item(96, 9).Value = 26060
$ws.Cells.Item(96, 12).Value = 0.05

# Row 97
$ws.Cells.Item(97, 3).Value = 7000
$ws.Cells.Item(97, 8).Value = 10832
$ws.Cells.Item(97, 9).Value = 3832
$ws.Cells.Item(97, 12).Value = 0.05

# Row 98
$ws.Cells.Item(98, 3).Value = 20000
$ws.Cells.Item(98, 8).Value = 143488
$ws.Cells.Item(98, 9).Value = 123488
$ws.Cells.Item(98, 12).Value = 0.05

# Row 99
$ws.Cells.Item(99, 3).Value = 20706
$ws.Cells.Item(99, 8).Value = 67869
$ws.Cells.Item(99, 9).Value = 47163
$ws.Cells.Item(99, 12).Value = 0.05

# Row 100
$ws.Cells.Item(100, 3).Value = 58000
$ws.Cells.Item(100, 8).Value = 175758
$ws.Cells.Item(100, 9).Value = 117758
$ws.Cells.Item(100, 12).Value = 0.05

# Row 101
$ws.Cells.Item(101, 3).Value = 36700
$ws.Cells.Item(101, 8).Value = 47662
$ws.Cells.Item(101, 9).Value = 10962
$ws.Cells.Item(101, 12).Value = 0.05

# Row 102
$ws.Cells.Item(102, 3).Value = 31000
$ws.Cells.Item(102, 8).Value = 35227
$ws.Cells.Item(102, 9).Value = 4227
$ws.Cells.Item(102, 12).Value = 0.05

# Row 103
$ws.Cells.Item(103, 3).Value = 73400
$ws.Cells.Item(103, 8).Value = 79532
$ws.Cells.Item(103, 9).Value = 6132
$ws.Cells.Item(103, 12).Value = 0.05

# Row 104
$ws.Cells.Item(104, 3).Value = 18500
$ws.Cells.Item(104, 8).Value = 18949
$ws.Cells.Item(104, 9).Value = 449
$ws.Cells.Item(104, 12).Value = 0.05

# Row 105
$ws.Cells.Item(105, 3).Value = 3900
$ws.Cells.Item(105, 8).Value = 3930
$ws.Cells.Item(105, 9).Value = 30
$ws.Cells.Item(105, 12).Value = 0.05

# Row 106
$ws.Cells.Item(106, 3).Value = 30800
$ws.Cells.Item(106, 8).Value = 39677
$ws.Cells.Item(106, 9).Value = 8877
$ws.Cells.Item(106, 12).Value = 0.05

# Row 107
$ws.Cells.Item(107, 3).Value = 40300
$ws.Cells.Item(107, 8).Value = 42094
$ws.Cells.Item(107, 9).Value = 1794
$ws.Cells.Item(107, 12).Value = 0.05

# Row 108
$ws.Cells.Item(108, 3).Value = 40600
$ws.Cells.Item(108, 8).Value = 40656
$ws.Cells.Item(108, 9).Value = 56
$ws.Cells.Item(108, 12).Value = 0.05

# Row 109
$ws.Cells.Item(109, 3).Value = 31400
$ws.Cells.Item(109, 8).Value = 32201
$ws.Cells.Item(109, 9).Value = 801
$ws.Cells.Item(109, 12).Value = 0.05

# Row 110
$ws.Cells.Item(110, 3).Value = 38100
$ws.Cells.Item(110, 8).Value = 49295
$ws.Cells.Item(110, 9).Value = 11195
$ws.Cells.Item(110, 12).Value = 0.05

# Row 111
$ws.Cells.Item(111, 3).Value = 27700
$ws.Cells.Item(111, 8).Value = 36772
$ws.Cells.Item(111, 9).Value = 9072
$ws.Cells.Item(111, 12).Value = 0.05

# Row 112
$ws.Cells.Item(112, 3).Value = 180500
$ws.Cells.Item(112, 8).Value = 217239
$ws.Cells.Item(112, 9).Value = 36739
$ws.Cells.Item(112, 12).Value = 0.05

# Row 113
$ws.Cells.Item(113, 3).Value = 17400
$ws.Cells.Item(113, 8).Value = 52598
$ws.Cells.Item(113, 9).Value = 35198
$ws.Cells.Item(113, 12).Value = 0.05

# Row 114
$ws.Cells.Item(114, 3).Value = 4400
$ws.Cells.Item(114, 8).Value = 4489
$ws.Cells.Item(114, 9).Value = 89
$ws.Cells.Item(114, 12).Value = 0.05

# Row 115
$ws.Cells.Item(115, 3).Value = 59900
$ws.Cells.Item(115, 8).Value = 61811
$ws.Cells.Item(115, 9).Value = 1911
$ws.Cells.Item(115, 12).Value = 0.05

# Row 116
$ws.Cells.Item(116, 3).Value = 46200
$ws.Cells.Item(116, 8).Value = 59613
$ws.Cells.Item(116, 9).Value = 13413
$ws.Cells.Item(116, 12).Value = 0.05

# Row 117
$ws.Cells.Item(117, 3).Value = 92100
$ws.Cells.Item(117, 8).Value = 119301
$ws.Cells.Item(117, 9).Value = 27201
$ws.Cells.Item(117, 12).Value = 0.05

# Row 118
$ws.Cells.Item(118, 3).Value = 13400
$ws.Cells.Item(118, 8).Value = 14709
$ws.Cells.Item(118, 9).Value = 1309
$ws.Cells.Item(118, 12).Value = 0.05

# Row 119
$ws.Cells.Item(119, 3).Value = 25100
$ws.Cells.Item(119, 8).Value = 25354
$ws.Cells.Item(119, 9).Value = 254
$ws.Cells.Item(119, 12).Value = 0.05

# Row 120
$ws.Cells.Item(120, 3).Value = 19900
$ws.Cells.Item(120, 8).Value = 21940
$ws.Cells.Item(120, 9).Value = 2040
$ws.Cells.Item(120, 12).Value = 0.05

# Row 121
$ws.Cells.Item(121, 3).Value = 17700
$ws.Cells.Item(121, 8).Value = 21612
$ws.Cells.Item(121, 9).Value = 3912
$ws.Cells.Item(121, 12).Value = 0.05

# Row 122
$ws.Cells.Item(122, 3).Value = 3300
$ws.Cells.Item(122, 8).Value = 3537
$ws.Cells.Item(122, 9).Value = 237
$ws.Cells.Item(122, 12).Value = 0.05

# Row 123
$ws.Cells.Item(123, 3).Value = 2600
$ws.Cells.Item(123, 8).Value = 2610
$ws.Cells.Item(123, 9).Value = 10
$ws.Cells.Item(123, 12).Value = 0.05

# Row 124
$ws.Cells.Item(124, 3).Value = 1300
$ws.Cells.Item(124, 8).Value = 1601
$ws.Cells.Item(124, 9).Value = 301
$ws.Cells.Item(124, 12).Value = 0.05

# Row 125
$ws.Cells.Item(125, 3).Value = 3600
$ws.Cells.Item(125, 8).Value = 3897
$ws.Cells.Item(125, 9).Value = 297
$ws.Cells.Item(125, 12).Value = 0.05

# Row 126
$ws.Cells.Item(126, 3).Value = 12500
$ws.Cells.Item(126, 8).Value = 12781
$ws.Cells.Item(126, 9).Value = 281
$ws.Cells.Item(126, 12).Value = 0.05

# Row 127
$ws.Cells.Item(127, 3).Value = 13100
$ws.Cells.Item(127, 8).Value = 13100
$ws.Cells.Item(127, 12).Value = 0.399812753899036

# Row 128
$ws.Cells.Item(128, 8).Value = 36055.16288364011
$ws.Cells.Item(128, 9).Value = 6055.162883640109
$ws.Cells.Item(128, 12).Value = 0.399812753899036

# Row 129
$ws.Cells.Item(129, 8).Value = 56370.16995668155
$ws.Cells.Item(129, 9).Value = 26370.16995668155
$ws.Cells.Item(129, 12).Value = 0.399812753899036


# Row 143 (full replace)
$ws.Range("A143:M143").ClearContents()
$ws.Cells.Item(143, 1).Value = 1918
$ws.Cells.Item(143, 2).Value = "HED"
$ws.Cells.Item(143, 3).Value = 7000
$ws.Cells.Item(143, 10).Value = 0
$ws.Cells.Item(143, 12).Value = 0.399812753899036
$ws.Cells.Item(143, 13).Value = 0.2

# Row 144 (full replace)
$ws.Range("A144:M144").ClearContents()
$ws.Cells.Item(144, 1).Value = 1919
$ws.Cells.Item(144, 2).Value = "HED"
$ws.Cells.Item(144, 3).Value = 8000
$ws.Cells.Item(144, 10).Value = 0
$ws.Cells.Item(144, 12).Value = 0.399812753899036
$ws.Cells.Item(144, 13).Value = 0.2

# Row 145 (full replace)
$ws.Range("A145:M145").ClearContents()
$ws.Cells.Item(145, 1).Value = 1920
$ws.Cells.Item(145, 2).Value = "HED"
$ws.Cells.Item(145, 3).Value = 38000
$ws.Cells.Item(145, 10).Value = 0
$ws.Cells.Item(145, 12).Value = 0.399812753899036
$ws.Cells.Item(145, 13).Value = 0.2

# Row 146 (full replace)
$ws.Range("A146:M146").ClearContents()
$ws.Cells.Item(146, 1).Value = 1921
$ws.Cells.Item(146, 2).Value = "HED"
$ws.Cells.Item(146, 3).Value = 1000
$ws.Cells.Item(146, 10).Value = 0
$ws.Cells.Item(146, 12).Value = 0.399812753899036
$ws.Cells.Item(146, 13).Value = 0.2

# Row 147 (full replace)
$ws.Range("A147:M147").ClearContents()
$ws.Cells.Item(147, 1).Value = 1922
$ws.Cells.Item(147, 2).Value = "HED"
$ws.Cells.Item(147, 3).Value = 70000
$ws.Cells.Item(147, 10).Value = 0
$ws.Cells.Item(147, 12).Value = 0.399812753899036
$ws.Cells.Item(147, 13).Value = 0.2


# Row 148 (new)
$ws.Cells.Item(148, 1).Value = 1923
$ws.Cells.Item(148, 2).Value = "HED"
$ws.Cells.Item(148, 3).Value = 90000
$ws.Cells.Item(148, 10).Value = 0
$ws.Cells.Item(148, 12).Value = 0.399812753899036
$ws.Cells.Item(148, 13).Value = 0.2

# Row 149 (new)
$ws.Cells.Item(149, 1).Value = 1924
$ws.Cells.Item(149, 2).Value = "HED"
$ws.Cells.Item(149, 3).Value = 120000
$ws.Cells.Item(149, 10).Value = 0
$ws.Cells.Item(149, 12).Value = 0.399812753899036
$ws.Cells.Item(149, 13).Value = 0.2

# Row 150 (new)
$ws.Cells.Item(150, 1).Value = 1925
$ws.Cells.Item(150, 2).Value = "HED"
$ws.Cells.Item(150, 3).Value = 80000
$ws.Cells.Item(150, 10).Value = 0
$ws.Cells.Item(150, 12).Value = 0.399812753899036
$ws.Cells.Item(150, 13).Value = 0.2

# Row 151 (new)
$ws.Cells.Item(151, 1).Value = 1926
$ws.Cells.Item(151, 2).Value = "HED"
$ws.Cells.Item(151, 3).Value = 65000
$ws.Cells.Item(151, 10).Value = 0
$ws.Cells.Item(151, 12).Value = 0.399812753899036
$ws.Cells.Item(151, 13).Value = 0.2

# Row 152 (new)
$ws.Cells.Item(152, 1).Value = 1927
$ws.Cells.Item(152, 2).Value = "HED"
$ws.Cells.Item(152, 3).Value = 70000
$ws.Cells.Item(152, 10).Value = 0
$ws.Cells.Item(152, 12).Value = 0.399812753899036
$ws.Cells.Item(152, 13).Value = 0.2

# Row 153 (new)
$ws.Cells.Item(153, 1).Value = 1928
$ws.Cells.Item(153, 2).Value = "HED"
$ws.Cells.Item(153, 3).Value = 70000
$ws.Cells.Item(153, 10).Value = 0
$ws.Cells.Item(153, 12).Value = 0.399812753899036
$ws.Cells.Item(153, 13).Value = 0.2

# Row 154 (new)
$ws.Cells.Item(154, 1).Value = 1929
$ws.Cells.Item(154, 2).Value = "HED"
$ws.Cells.Item(154, 3).Value = 135000
$ws.Cells.Item(154, 10).Value = 0
$ws.Cells.Item(154, 12).Value = 0.399812753899036
$ws.Cells.Item(154, 13).Value = 0.2

# Row 155 (new)
$ws.Cells.Item(155, 1).Value = 1930
$ws.Cells.Item(155, 2).Value = "HED"
$ws.Cells.Item(155, 3).Value = 40000
$ws.Cells.Item(155, 10).Value = 0
$ws.Cells.Item(155, 12).Value = 0.399812753899036
$ws.Cells.Item(155, 13).Value = 0.2

# Row 156 (new)
$ws.Cells.Item(156, 1).Value = 1931
$ws.Cells.Item(156, 2).Value = "HED"
$ws.Cells.Item(156, 3).Value = 50000
$ws.Cells.Item(156, 10).Value = 0
$ws.Cells.Item(156, 12).Value = 0.399812753899036
$ws.Cells.Item(156, 13).Value = 0.2

# Row 157 (new)
$ws.Cells.Item(157, 1).Value = 1932
$ws.Cells.Item(157, 2).Value = "HED"
$ws.Cells.Item(157, 3).Value = 35000
$ws.Cells.Item(157, 10).Value = 0
$ws.Cells.Item(157, 12).Value = 0.399812753899036
$ws.Cells.Item(157, 13).Value = 0.2

# Row 158 (new)
$ws.Cells.Item(158, 1).Value = 1933
$ws.Cells.Item(158, 2).Value = "HED"
$ws.Cells.Item(158, 3).Value = 7500
$ws.Cells.Item(158, 10).Value = 0
$ws.Cells.Item(158, 12).Value = 0.399812753899036
$ws.Cells.Item(158, 13).Value = 0.2

# Row 159 (new)
$ws.Cells.Item(159, 1).Value = 1934
$ws.Cells.Item(159, 2).Value = "HED"
$ws.Cells.Item(159, 3).Value = 15000
$ws.Cells.Item(159, 10).Value = 0
$ws.Cells.Item(159, 12).Value = 0.399812753899036
$ws.Cells.Item(159, 13).Value = 0.2

# Row 160 (new)
$ws.Cells.Item(160, 1).Value = 1935
$ws.Cells.Item(160, 2).Value = "HED"
$ws.Cells.Item(160, 3).Value = 45000
$ws.Cells.Item(160, 10).Value = 0
$ws.Cells.Item(160, 12).Value = 0.399812753899036
$ws.Cells.Item(160, 13).Value = 0.2

# Row 161 (new)
$ws.Cells.Item(161, 1).Value = 1936
$ws.Cells.Item(161, 2).Value = "HED"
$ws.Cells.Item(161, 3).Value = 2000
$ws.Cells.Item(161, 10).Value = 0
$ws.Cells.Item(161, 12).Value = 0.399812753899036
$ws.Cells.Item(161, 13).Value = 0.2

# Row 162 (new)
$ws.Cells.Item(162, 1).Value = 1937
$ws.Cells.Item(162, 2).Value = "HED"
$ws.Cells.Item(162, 3).Value = 38000
$ws.Cells.Item(162, 10).Value = 0
$ws.Cells.Item(162, 12).Value = 0.399812753899036
$ws.Cells.Item(162, 13).Value = 0.2

# Row 163 (new)
$ws.Cells.Item(163, 1).Value = 1938
$ws.Cells.Item(163, 2).Value = "HED"
$ws.Cells.Item(163, 3).Value = 10000
$ws.Cells.Item(163, 10).Value = 0
$ws.Cells.Item(163, 12).Value = 0.399812753899036
$ws.Cells.Item(163, 13).Value = 0.2

# Row 164 (new)
$ws.Cells.Item(164, 1).Value = 1939
$ws.Cells.Item(164, 2).Value = "HED"
$ws.Cells.Item(164, 3).Value = 14000
$ws.Cells.Item(164, 10).Value = 0
$ws.Cells.Item(164, 12).Value = 0.399812753899036
$ws.Cells.Item(164, 13).Value = 0.2

# Row 165 (new)
$ws.Cells.Item(165, 1).Value = 1940
$ws.Cells.Item(165, 2).Value = "HED"
$ws.Cells.Item(165, 3).Value = 60000
$ws.Cells.Item(165, 10).Value = 0
$ws.Cells.Item(165, 12).Value = 0.399812753899036
$ws.Cells.Item(165, 13).Value = 0.2

# Row 166 (new)
$ws.Cells.Item(166, 1).Value = 1941
$ws.Cells.Item(166, 2).Value = "HED"
$ws.Cells.Item(166, 3).Value = 2000
$ws.Cells.Item(166, 10).Value = 0
$ws.Cells.Item(166, 12).Value = 0.399812753899036
$ws.Cells.Item(166, 13).Value = 0.2

# Row 167 (new)
$ws.Cells.Item(167, 1).Value = 1942
$ws.Cells.Item(167, 2).Value = "HED"
$ws.Cells.Item(167, 3).Value = 7000
$ws.Cells.Item(167, 10).Value = 0
$ws.Cells.Item(167, 12).Value = 0.399812753899036
$ws.Cells.Item(167, 13).Value = 0.2

# Row 168 (new)
$ws.Cells.Item(168, 1).Value = 1943
$ws.Cells.Item(168, 2).Value = "HED"
$ws.Cells.Item(168, 3).Value = 5000
$ws.Cells.Item(168, 10).Value = 0
$ws.Cells.Item(168, 12).Value = 0.399812753899036
$ws.Cells.Item(168, 13).Value = 0.2

# Row 169 (new)
$ws.Cells.Item(169, 1).Value = 1944
$ws.Cells.Item(169, 2).Value = "HED"
$ws.Cells.Item(169, 3).Value = 5000
$ws.Cells.Item(169, 10).Value = 0
$ws.Cells.Item(169, 12).Value = 0.399812753899036
$ws.Cells.Item(169, 13).Value = 0.2

# Row 170 (new)
$ws.Cells.Item(170, 1).Value = 1945
$ws.Cells.Item(170, 2).Value = "HED"
$ws.Cells.Item(170, 3).Value = 14000
$ws.Cells.Item(170, 10).Value = 0
$ws.Cells.Item(170, 12).Value = 0.399812753899036
$ws.Cells.Item(170, 13).Value = 0.2

# Row 171 (new)
$ws.Cells.Item(171, 1).Value = 1946
$ws.Cells.Item(171, 2).Value = "HED"
$ws.Cells.Item(171, 3).Value = 14000
$ws.Cells.Item(171, 10).Value = 0
$ws.Cells.Item(171, 12).Value = 0.399812753899036
$ws.Cells.Item(171, 13).Value = 0.2

# Row 172 (new)
$ws.Cells.Item(172, 1).Value = 1947
$ws.Cells.Item(172, 2).Value = "HED"
$ws.Cells.Item(172, 3).Value = 7000
$ws.Cells.Item(172, 10).Value = 0
$ws.Cells.Item(172, 12).Value = 0.399812753899036
$ws.Cells.Item(172, 13).Value = 0.2

# Row 173 (new)
$ws.Cells.Item(173, 1).Value = 1948
$ws.Cells.Item(173, 2).Value = "HED"
$ws.Cells.Item(173, 3).Value = 7000
$ws.Cells.Item(173, 10).Value = 0
$ws.Cells.Item(173, 12).Value = 0.399812753899036
$ws.Cells.Item(173, 13).Value = 0.2

# Row 174 (new)
$ws.Cells.Item(174, 1).Value = 1949
$ws.Cells.Item(174, 2).Value = "HED"
$ws.Cells.Item(174, 3).Value = 30000
$ws.Cells.Item(174, 10).Value = 0
$ws.Cells.Item(174, 12).Value = 0.399812753899036
$ws.Cells.Item(174, 13).Value = 0.2

# Row 175 (new)
$ws.Cells.Item(175, 1).Value = 1950
$ws.Cells.Item(175, 2).Value = "HED"
$ws.Cells.Item(175, 3).Value = 14000
$ws.Cells.Item(175, 10).Value = 0
$ws.Cells.Item(175, 12).Value = 0.399812753899036
$ws.Cells.Item(175, 13).Value = 0.2

# Row 176 (new)
$ws.Cells.Item(176, 1).Value = 1951
$ws.Cells.Item(176, 2).Value = "HED"
$ws.Cells.Item(176, 3).Value = 14000
$ws.Cells.Item(176, 10).Value = 0
$ws.Cells.Item(176, 12).Value = 0.399812753899036
$ws.Cells.Item(176, 13).Value = 0.2

# Row 177 (new)
$ws.Cells.Item(177, 1).Value = 1952
$ws.Cells.Item(177, 2).Value = "HED"
$ws.Cells.Item(177, 3).Value = 14000
$ws.Cells.Item(177, 10).Value = 0
$ws.Cells.Item(177, 12).Value = 0.399812753899036
$ws.Cells.Item(177, 13).Value = 0.2

# Row 178 (new)
$ws.Cells.Item(178, 1).Value = 1953
$ws.Cells.Item(178, 2).Value = "HED"
$ws.Cells.Item(178, 3).Value = 14000
$ws.Cells.Item(178, 10).Value = 0
$ws.Cells.Item(178, 12).Value = 0.399812753899036
$ws.Cells.Item(178, 13).Value = 0.2

# Row 179 (new)
$ws.Cells.Item(179, 1).Value = 1954
$ws.Cells.Item(179, 2).Value = "HED"
$ws.Cells.Item(179, 3).Value = 30000
$ws.Cells.Item(179, 10).Value = 0
$ws.Cells.Item(179, 12).Value = 0.399812753899036
$ws.Cells.Item(179, 13).Value = 0.2

# Row 180 (new)
$ws.Cells.Item(180, 1).Value = 1955
$ws.Cells.Item(180, 2).Value = "HED"
$ws.Cells.Item(180, 3).Value = 14000
$ws.Cells.Item(180, 10).Value = 0
$ws.Cells.Item(180, 12).Value = 0.399812753899036
$ws.Cells.Item(180, 13).Value = 0.2

# Row 181 (new)
$ws.Cells.Item(181, 1).Value = 1956
$ws.Cells.Item(181, 2).Value = "HED"
$ws.Cells.Item(181, 3).Value = 7000
$ws.Cells.Item(181, 10).Value = 0
$ws.Cells.Item(181, 12).Value = 0.399812753899036
$ws.Cells.Item(181, 13).Value = 0.2

# Row 182 (new)
$ws.Cells.Item(182, 1).Value = 1957
$ws.Cells.Item(182, 2).Value = "HED"
$ws.Cells.Item(182, 3).Value = 7000
$ws.Cells.Item(182, 10).Value = 0
$ws.Cells.Item(182, 12).Value = 0.399812753899036
$ws.Cells.Item(182, 13).Value = 0.2

# Row 183 (new)
$ws.Cells.Item(183, 1).Value = 1958
$ws.Cells.Item(183, 2).Value = "HED"
$ws.Cells.Item(183, 3).Value = 14000
$ws.Cells.Item(183, 10).Value = 0
$ws.Cells.Item(183, 12).Value = 0.399812753899036
$ws.Cells.Item(183, 13).Value = 0.2

# Row 184 (new)
$ws.Cells.Item(184, 1).Value = 1959
$ws.Cells.Item(184, 2).Value = "HED"
$ws.Cells.Item(184, 3).Value = 14000
$ws.Cells.Item(184, 10).Value = 0
$ws.Cells.Item(184, 12).Value = 0.399812753899036
$ws.Cells.Item(184, 13).Value = 0.2

# Row 185 (new)
$ws.Cells.Item(185, 1).Value = 1960
$ws.Cells.Item(185, 2).Value = "HED"
$ws.Cells.Item(185, 3).Value = 10000
$ws.Cells.Item(185, 10).Value = 0
$ws.Cells.Item(185, 12).Value = 0.399812753899036
$ws.Cells.Item(185, 13).Value = 0.2

# Row 186 (new)
$ws.Cells.Item(186, 1).Value = 1961
$ws.Cells.Item(186, 2).Value = "HED"
$ws.Cells.Item(186, 3).Value = 14000
$ws.Cells.Item(186, 10).Value = 0
$ws.Cells.Item(186, 12).Value = 0.399812753899036
$ws.Cells.Item(186, 13).Value = 0.2

# Row 187 (new)
$ws.Cells.Item(187, 1).Value = 1962
$ws.Cells.Item(187, 2).Value = "HED"
$ws.Cells.Item(187, 3).Value = 30000
$ws.Cells.Item(187, 10).Value = 0
$ws.Cells.Item(187, 12).Value = 0.399812753899036
$ws.Cells.Item(187, 13).Value = 0.2

# Row 188 (new)
$ws.Cells.Item(188, 1).Value = 1963
$ws.Cells.Item(188, 2).Value = "HED"
$ws.Cells.Item(188, 3).Value = 18000
$ws.Cells.Item(188, 10).Value = 0
$ws.Cells.Item(188, 12).Value = 0.399812753899036
$ws.Cells.Item(188, 13).Value = 0.2

# Row 189 (new)
$ws.Cells.Item(189, 1).Value = 1964
$ws.Cells.Item(189, 2).Value = "HED"
$ws.Cells.Item(189, 3).Value = 30000
$ws.Cells.Item(189, 10).Value = 0
$ws.Cells.Item(189, 12).Value = 0.399812753899036
$ws.Cells.Item(189, 13).Value = 0.2

# Row 190 (new)
$ws.Cells.Item(190, 1).Value = 1965
$ws.Cells.Item(190, 2).Value = "HED"
$ws.Cells.Item(190, 3).Value = 18000
$ws.Cells.Item(190, 10).Value = 0
$ws.Cells.Item(190, 12).Value = 0.399812753899036
$ws.Cells.Item(190, 13).Value = 0.2

# Row 191 (new)
$ws.Cells.Item(191, 1).Value = 1966
$ws.Cells.Item(191, 2).Value = "HED"
$ws.Cells.Item(191, 3).Value = 80000
$ws.Cells.Item(191, 10).Value = 0
$ws.Cells.Item(191, 12).Value = 0.399812753899036
$ws.Cells.Item(191, 13).Value = 0.2

# Row 192 (new)
$ws.Cells.Item(192, 1).Value = 1967
$ws.Cells.Item(192, 2).Value = "HED"
$ws.Cells.Item(192, 3).Value = 80000
$ws.Cells.Item(192, 10).Value = 0
$ws.Cells.Item(192, 12).Value = 0.399812753899036
$ws.Cells.Item(192, 13).Value = 0.2

# Row 193 (new)
$ws.Cells.Item(193, 1).Value = 1968
$ws.Cells.Item(193, 2).Value = "HED"
$ws.Cells.Item(193, 3).Value = 30000
$ws.Cells.Item(193, 10).Value = 0
$ws.Cells.Item(193, 12).Value = 0.399812753899036
$ws.Cells.Item(193, 13).Value = 0.2

# Row 194 (new)
$ws.Cells.Item(194, 1).Value = 1969
$ws.Cells.Item(194, 2).Value = "HED"
$ws.Cells.Item(194, 3).Value = 18000
$ws.Cells.Item(194, 10).Value = 0
$ws.Cells.Item(194, 12).Value = 0.399812753899036
$ws.Cells.Item(194, 13).Value = 0.2

# Row 195 (new)
$ws.Cells.Item(195, 1).Value = 1970
$ws.Cells.Item(195, 2).Value = "HED"
$ws.Cells.Item(195, 3).Value = 9000
$ws.Cells.Item(195, 8).Value = 14317
$ws.Cells.Item(195, 9).Value = 5317
$ws.Cells.Item(195, 10).Value = 0
$ws.Cells.Item(195, 12).Value = 0.05
$ws.Cells.Item(195, 13).Value = 0.2

# Row 196 (new)
$ws.Cells.Item(196, 1).Value = 1971
$ws.Cells.Item(196, 2).Value = "HED"
$ws.Cells.Item(196, 3).Value = 7500
$ws.Cells.Item(196, 8).Value = 8522
$ws.Cells.Item(196, 9).Value = 1022
$ws.Cells.Item(196, 10).Value = 0
$ws.Cells.Item(196, 12).Value = 0.05
$ws.Cells.Item(196, 13).Value = 0.2

# Row 197 (new)
$ws.Cells.Item(197, 1).Value = 1972
$ws.Cells.Item(197, 2).Value = "HED"
$ws.Cells.Item(197, 3).Value = 3500
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 0.91
$ws.Cells.Item(197, 6).Value = 0.06999999999999999
$ws.Cells.Item(197, 7).Value = 0.02
$ws.Cells.Item(197, 8).Value = 4407
$ws.Cells.Item(197, 9).Value = 907
$ws.Cells.Item(197, 10).Value = 43
$ws.Cells.Item(197, 12).Value = 0.05
$ws.Cells.Item(197, 13).Value = 0.2

# Row 198 (new)
$ws.Cells.Item(198, 1).Value = 1973
$ws.Cells.Item(198, 2).Value = "HED"
$ws.Cells.Item(198, 3).Value = 40000
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 5).Value = 0.66
$ws.Cells.Item(198, 6).Value = 0.34
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 58495
$ws.Cells.Item(198, 9).Value = 18495
$ws.Cells.Item(198, 10).Value = 99
$ws.Cells.Item(198, 12).Value = 0.05
$ws.Cells.Item(198, 13).Value = 0.2

# Row 199 (new)
$ws.Cells.Item(199, 1).Value = 1974
$ws.Cells.Item(199, 2).Value = "HED"
$ws.Cells.Item(199, 3).Value = 6000
$ws.Cells.Item(199, 4).Value = 0
$ws.Cells.Item(199, 5).Value = 0.57
$ws.Cells.Item(199, 6).Value = 0.43
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 8936
$ws.Cells.Item(199, 9).Value = 2936
$ws.Cells.Item(199, 10).Value = 113
$ws.Cells.Item(199, 12).Value = 0.05
$ws.Cells.Item(199, 13).Value = 0.2

# Row 200 (new)
$ws.Cells.Item(200, 1).Value = 1975
$ws.Cells.Item(200, 2).Value = "HED"
$ws.Cells.Item(200, 3).Value = 10000
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 0.5742574257425742
$ws.Cells.Item(200, 6).Value = 0.4158415841584159
$ws.Cells.Item(200, 7).Value = 0.009900990099009901
$ws.Cells.Item(200, 8).Value = 17635
$ws.Cells.Item(200, 9).Value = 7635
$ws.Cells.Item(200, 10).Value = 195
$ws.Cells.Item(200, 12).Value = 0.05
$ws.Cells.Item(200, 13).Value = 0.2

# Row 201 (new)
$ws.Cells.Item(201, 1).Value = 1976
$ws.Cells.Item(201, 2).Value = "HED"
$ws.Cells.Item(201, 3).Value = 3500
$ws.Cells.Item(201, 4).Value = 0
$ws.Cells.Item(201, 5).Value = 0.7319587628865979
$ws.Cells.Item(201, 6).Value = 0.2680412371134021
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 8825
$ws.Cells.Item(201, 9).Value = 5325
$ws.Cells.Item(201, 10).Value = 191
$ws.Cells.Item(201, 12).Value = 0.05
$ws.Cells.Item(201, 13).Value = 0.2

